$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 295
$ws.Range("I18").Value = 295
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 295
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -11

# Row 106
$ws.Range("H106").Value = 47621508
$ws.Range("I106").Value = 58825960
$ws.Range("J106").Value = 2597.5
$ws.Range("K106").Value = 58825960
$ws.Range("L106").Value = 2597.5
$ws.Range("M106").Value = -58825329
$ws.Range("N106").Value = -3859.5

# Row 135
$ws.Range("H135").Value = 5556210
$ws.Range("I135").Value = 5814544.5
$ws.Range("J135").Value = 2018
$ws.Range("K135").Value = 52330900.5
$ws.Range("L135").Value = 18162
$ws.Range("M135").Value = -52328365.5
$ws.Range("N135").Value = -23232

# Row 137
$ws.Range("H137").Value = 2320.6
$ws.Range("I137").Value = 1050
$ws.Range("J137").Value = 3167.6667
$ws.Range("K137").Value = 3150
$ws.Range("L137").Value = 9503.000100000001
$ws.Range("N137").Value = -14603.0001
$ws.Range("M137").Value = -600

# Row 141
$ws.Range("H141").Value = 4672.8423
$ws.Range("I141").Value = 929.44446
$ws.Range("J141").Value = 8041.9
$ws.Range("K141").Value = 2788.33338
$ws.Range("L141").Value = 24125.7
$ws.Range("M141").Value = 2391.66662
$ws.Range("N141").Value = -34485.7

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1267.25
$ws.Range("I74").Value = 1458.1428
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1458.1428
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -584.1428000000001
$ws.Range("N74").Value = -2748

# Row 77
$ws.Range("H77").Value = 1267.25
$ws.Range("I77").Value = 1458.1428
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 7290.714
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -2922.714
$ws.Range("N77").Value = -13736

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4832.0967
$ws.Range("I134").Value = 1546.25
$ws.Range("J134").Value = 35500
$ws.Range("K134").Value = 4638.75
$ws.Range("L134").Value = 106500
$ws.Range("M134").Value = -2103.75
$ws.Range("N134").Value = -111570

# Row 139
$ws.Range("H139").Value = 40859.8
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 40859.8
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 40859.8
$ws.Range("N139").Value = -51139.8

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""

# Row 7
$ws.Range("H7").Value = 7286.5
$ws.Range("I7").Value = 167.625
$ws.Range("J7").Value = 16778.334
$ws.Range("K7").Value = 167.625
$ws.Range("L7").Value = 16778.334
$ws.Range("M7").Value = -54.625
$ws.Range("N7").Value = -17004.334

# Row 17
$ws.Range("H17").Value = 20280
$ws.Range("I17").Value = 990
$ws.Range("J17").Value = 29925
$ws.Range("K17").Value = 990
$ws.Range("L17").Value = 29925
$ws.Range("M17").Value = -816
$ws.Range("N17").Value = -30273

# Row 25
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = -10348

# Row 41
$ws.Range("H41").Value = 3779.5
$ws.Range("I41").Value = 3779.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3779.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -3351.5

# Row 53
$ws.Range("H53").Value = 40500
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 40500
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 40500
$ws.Range("N53").Value = -41714

# Row 58
$ws.Range("H58").Value = 43479050
$ws.Range("I58").Value = 71429340
$ws.Range("J58").Value = 811.44446
$ws.Range("K58").Value = 71429340
$ws.Range("L58").Value = 811.44446
$ws.Range("M58").Value = -71429137
$ws.Range("N58").Value = -1217.44446

# Row 59
$ws.Range("H59").Value = 40000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 40000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 40000
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -42290

# Row 60
$ws.Range("H60").Value = 11275.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 11275.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11275.5
$ws.Range("N60").Value = -12297.5

# Row 68
$ws.Range("H68").Value = 19757.4
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 19757.4
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 19757.4
$ws.Range("N68").Value = -21255.4

# Row 71
$ws.Range("H71").Value = 19757.4
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 19757.4
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 59272.2
$ws.Range("N71").Value = -66760.20000000001

# Row 74
$ws.Range("H74").Value = 18316
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 16645
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 16645
$ws.Range("M74").Value = -24126
$ws.Range("N74").Value = -18393

# Row 77
$ws.Range("H77").Value = 18316
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 16645
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 49935
$ws.Range("M77").Value = -70632
$ws.Range("N77").Value = -58671

# Row 98
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""

# Row 107
$ws.Range("H107").Value = 678.375
$ws.Range("I107").Value = 539
$ws.Range("J107").Value = 817.75
$ws.Range("K107").Value = 539
$ws.Range("L107").Value = 817.75
$ws.Range("M107").Value = 1381
$ws.Range("N107").Value = -4657.75

# Row 136
$ws.Range("H136").Value = 43479050
$ws.Range("I136").Value = 71429340
$ws.Range("J136").Value = 811.44446
$ws.Range("K136").Value = 214288020
$ws.Range("L136").Value = 2434.33338
$ws.Range("M136").Value = -214285470
$ws.Range("N136").Value = -7534.33338

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 13862.75
$ws.Range("I33").Value = 1750.25
$ws.Range("J33").Value = 25975.25
$ws.Range("K33").Value = 10501.5
$ws.Range("L33").Value = 155851.5
$ws.Range("M33").Value = -10218.5
$ws.Range("N33").Value = -156417.5

# Row 34
$ws.Range("H34").Value = 381.6154
$ws.Range("I34").Value = 359.5
$ws.Range("J34").Value = 417
$ws.Range("K34").Value = 1078.5
$ws.Range("L34").Value = 1251
$ws.Range("M34").Value = -994.5
$ws.Range("N34").Value = -1419

# Row 39
$ws.Range("H39").Value = 900
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 900
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 2700
$ws.Range("N39").Value = -3288

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5557818
$ws.Range("I126").Value = 1753.3334
$ws.Range("J126").Value = 6669031
$ws.Range("K126").Value = 5260.0002
$ws.Range("L126").Value = 20007093
$ws.Range("M126").Value = -2790.0002
$ws.Range("N126").Value = -20012033

# Row 135
$ws.Range("H135").Value = 46898.184
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 46898.184
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 46898.184
$ws.Range("N135").Value = -57038.184

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 36875000
$ws.Range("I132").Value = 95239660
$ws.Range("J132").Value = 13105.158
$ws.Range("K132").Value = 285718980
$ws.Range("L132").Value = 39315.474
$ws.Range("M132").Value = -285716450
$ws.Range("N132").Value = -44375.474
